$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Step 1: copy row 29 formatting down into new rows 30 and 31 ----
$ws.Range("A29:T29").Copy()
$ws.Range("A30:T30").PasteSpecial(-4122)
$ws.Range("A31:T31").PasteSpecial(-4122)

# ---- Step 2: row label (column B) updates for rows 4-31 ----
$labels = @{
    4 = "Holden"
    5 = "Rizzie Spiral"
    6 = "RotRing OmegaMax-90"
    7 = "Equal Angle"
    8 = "Tilt Rotate"
    9 = "CLR"
    10 = "Rizzie Hex"
    11 = "Matthies Hex"
    12 = "Tilt Rotate_Partial"
    13 = "RotRing OmegaMax-60"
    14 = "Equal Angle_Partial"
    15 = "Rizzie Hex_Partial"
    16 = "ND Single"
    17 = "RD Single"
    18 = "TD Single"
    19 = "Morris Single"
    20 = "Ring Perpendicular to ND"
    21 = "Ring Perpendicular to RD"
    22 = "Ring Perpendicular to TD"
    23 = "OffsetFTD"
    24 = "OffsetATD"
    25 = "OffsetF45"
    26 = "OffsetA45"
    27 = "OffsetFRD"
    28 = "OffsetARD"
    29 = "Gaussian Quadrature"
    30 = "Michael-CCHex"
    31 = "Michael-SNHex"
}
foreach ($r in $labels.Keys) {
    $ws.Cells.Item($r, 2).Value = $labels[$r]
}

# ---- Step 3: column A index values for new rows 30-31 ----
$ws.Cells.Item(30, 1).Value = 28
$ws.Cells.Item(31, 1).Value = 29

# ---- Step 4: numeric data C:T for rows 4-31 ----
$rowValues = @{
    4 = @(0.9935086146764005, 1.001698563082111, 1.001580782760231, 0.9935086146764005, 1.001615119672724, 1.001698563082111, 0.9964569955499462, 1.004192984755004, 1.001698563082111, 1.001580782760231, 0.997544698718316, 0.997544698718316, 0.997182130995526, 0.9989293201729144, 0.9989293201729144, 0.9996216309002135, 0.9996216309002135, 0.9998421767494028)
    5 = @(1.075336831766261, 0.9812902494817299, 0.981096656893965, 1.075336831766261, 0.9811531017487453, 0.9812902494817299, 1.041480059195947, 0.9495542042648314, 0.9812902494817299, 0.981096656893965, 1.028216744330113, 1.028216744330113, 1.032637849285391, 1.012574579380652, 1.012574579380652, 1.004753496905922, 1.004753496905922, 1.00165185055858)
    6 = @(1.024208985700529, 0.9938988017437617, 0.9939749498659368, 1.024208985700529, 0.9939527494518351, 0.9938988017437617, 1.013297318178215, 0.9839477027363513, 0.9938988017437617, 0.9939749498659368, 1.009091967783233, 1.009091967783233, 1.010493751248227, 1.004027579103409, 1.004027579103409, 1.001495384763497, 1.001495384763497, 1.000546751279438)
    7 = @(0.992969347233429, 1.001847905273776, 1.001707530518732, 0.992969347233429, 1.00174845371758, 1.001847905273776, 0.9961656277665722, 1.004526674697407, 1.001847905273776, 1.001707530518732, 0.9973384388760804, 0.9973384388760804, 0.9969475018395778, 0.9988415943419788, 0.9988415943419788, 0.999593172074928, 0.999593172074928, 0.9998275898679158)
    8 = @(0.9730148253937378, 1.007184783031468, 1.006502694454478, 0.9730148253937378, 1.00670153736151, 1.007184783031468, 0.985316014044316, 1.017210591230337, 1.007184783031468, 1.006502694454478, 0.9897587599241078, 0.9897587599241078, 0.9882778446308439, 0.9955674342932278, 0.9955674342932278, 0.9984717714777878, 0.9984717714777878, 0.9993217409193078)
    9 = @(0.9994992730491395, 1.00011613681466, 1.000130207142011, 0.9994992730491395, 1.000126105247291, 1.00011613681466, 0.9997213436553589, 1.000349900401669, 1.00011613681466, 1.000130207142011, 0.9998147400955755, 0.9998147400955755, 0.9997836079488366, 0.9999152056686039, 0.9999152056686039, 0.999965438455118, 0.999965438455118, 0.9999904943850217)
    10 = @(0.9999315057285559, 1.000019946708355, 1.000015556391327, 0.9999315057285559, 1.000016835909944, 1.000019946708355, 0.9999633441660268, 1.000040645159586, 1.000019946708355, 1.000015556391327, 0.9999735310599414, 0.9999735310599414, 0.9999701354286366, 0.999989002942746, 0.999989002942746, 0.9999967388841483, 0.9999967388841483, 0.9999979723439658)
    11 = @(0.9985466681248606, 1.000344976798933, 1.000373531000567, 0.9985466681248606, 1.000365205458988, 1.000344976798933, 0.999194060057356, 1.001001521959724, 1.000344976798933, 1.000373531000567, 0.9994600995627136, 0.9994600995627136, 0.9993714197275945, 0.9997550586414533, 0.9997550586414533, 0.9999025381808231, 0.9999025381808231, 0.9999709939000714)
    12 = @(0.9726070967102842, 1.007300375174958, 1.006597038371994, 0.9726070967102842, 1.006802075350439, 1.007300375174958, 0.9850966803889492, 1.017458127560284, 1.007300375174958, 1.006597038371994, 0.9896020675411392, 0.9896020675411392, 0.9881002718237425, 0.9955015034190787, 0.9955015034190787, 0.9984512213580484, 0.9984512213580484, 0.9993102322594848)
    13 = @(1.034541494822689, 0.9913611073668852, 0.991366581221883, 1.034541494822689, 0.9913649859854505, 0.9913611073668852, 1.018996553838803, 0.9769785926650824, 0.9913611073668852, 0.991366581221883, 1.012954038022286, 1.012954038022286, 1.014968209961125, 1.005756394470486, 1.005756394470486, 1.002157572694586, 1.002157572694586, 1.000768219316799)
    14 = @(0.9931202773157891, 1.001857266368422, 1.001643635600001, 0.9931202773157891, 1.00170591576842, 1.001857266368422, 0.9962655923473693, 1.004342334526315, 1.001857266368422, 1.001643635600001, 0.9973819564578952, 0.9973819564578952, 0.9970098350877198, 0.9988737264280708, 0.9988737264280708, 0.9996196114131586, 0.9996196114131586, 0.9998225036543862)
    15 = @(1.00533638118444, 0.9986472728469137, 0.9986762579753817, 1.00533638118444, 0.9986678058666526, 0.9986472728469137, 1.002928300764814, 0.9964755400407521, 0.9986472728469137, 0.9986762579753817, 1.002006319579911, 1.002006319579911, 1.002313646641545, 1.000886637335579, 1.000886637335579, 1.000326796213412, 1.000326796213412, 1.000121926446492)
    16 = @(0.9526597999999998, 1.012634, 1.0113912, 0.9526597999999998, 1.011753499999999, 1.012634, 0.9742505199999987, 1.0301398, 1.012634, 1.0113912, 0.9820254999999998, 0.9820254999999998, 0.9794338399999994, 0.9922283333333329, 0.9922283333333329, 0.9973297499999996, 0.9973297499999996, 0.9988048033333327)
    17 = @(0.9526598, 1.012634, 1.0113912, 0.9526598, 1.0117535, 1.012634, 0.9742505199999999, 1.0301398, 1.012634, 1.0113912, 0.9820255, 0.9820255, 0.97943384, 0.9922283333333333, 0.9922283333333333, 0.99732975, 0.99732975, 0.9988048033333333)
    18 = @(1.1873941, 0.9526598, 0.95342462, 1.1873941, 0.95320164, 0.9526598, 1.1028898, 0.8759446600000002, 0.9526598, 0.95342462, 1.07040936, 1.07040936, 1.081236173333333, 1.031159506666667, 1.031159506666667, 1.01153458, 1.01153458, 1.004252436666667)
    19 = @(0.9971564699999999, 1.0006024, 1.0007712, 0.9971564699999999, 1.000722, 1.0006024, 0.998397, 1.0020886, 1.0006024, 1.0007712, 0.9989638349999999, 0.9989638349999999, 0.99877489, 0.9995100233333333, 0.9995100233333333, 0.9997831174999999, 0.9997831174999999, 0.9999562783333333)
    20 = @(0.9822361227397259, 1.004205658630136, 1.004571701917808, 0.9822361227397259, 1.004464977260274, 1.004205658630136, 0.990145156712329, 1.012260926027397, 1.004205658630136, 1.004571701917808, 0.993403912328767, 0.993403912328767, 0.992317660456621, 0.9970044944292235, 0.9970044944292235, 0.9988047854794517, 0.9988047854794517, 0.9996474238812785)
    21 = @(0.987245854736842, 1.002960087894737, 1.003315456315789, 0.987245854736842, 1.003211843684211, 1.002960087894737, 0.9929029810526314, 1.008908904210527, 1.002960087894737, 1.003315456315789, 0.9952806555263156, 0.9952806555263156, 0.9944880973684209, 0.9978404663157893, 0.9978404663157893, 0.9991203717105261, 0.9991203717105261, 0.9997575213157894)
    22 = @(1.063849743157895, 0.9842251800000001, 0.9839333421052632, 1.063849743157895, 0.9840184273684212, 0.9842251800000001, 1.035184897894737, 0.9570999831578948, 0.9842251800000001, 0.9839333421052632, 1.023891542631579, 1.023891542631579, 1.027655994385965, 1.010669421754386, 1.010669421754386, 1.004058361315789, 1.004058361315789, 1.001385262280702)
    23 = @(0.903179289747316, 1.023710879888909, 1.024479788244299, 0.903179289747316, 1.024255608649295, 1.023710879888909, 0.9465706614042574, 1.065425904135564, 1.023710879888909, 1.024479788244299, 0.9638295389958076, 0.9638295389958076, 0.9580765797986243, 0.9837899859601746, 0.9837899859601746, 0.9937702094423582, 0.9937702094423582, 0.997937022011607)
    24 = @(0.9737550851655731, 1.006768384060493, 1.006446150929852, 0.9737550851655731, 1.00654008781614, 1.006768384060493, 0.9856398703219604, 1.01712833877372, 1.006768384060493, 1.006446150929852, 0.9901006180477125, 0.9901006180477125, 0.9886137021391285, 0.9956565400519727, 0.9956565400519728, 0.9984345010541029, 0.9984345010541029, 0.9993796528446234)
    25 = @(0.9987123041058847, 1.000296024782096, 1.000336325498456, 0.9987123041058847, 1.000324559131632, 1.000296024782096, 0.9992824357394359, 1.000904523457336, 1.000296024782096, 1.000336325498456, 0.9995243148021703, 0.9995243148021703, 0.9994436884479255, 0.9997815514621454, 0.9997815514621454, 0.9999101697921331, 0.9999101697921331, 0.9999760287858068)
    26 = @(0.9996364740418497, 1.000103377191377, 1.000083939203447, 0.9996364740418497, 1.000089609287539, 1.000103377191377, 0.9998045610426181, 1.000220137397193, 1.000103377191377, 1.000083939203447, 0.9998602066226483, 0.9998602066226483, 0.9998416580959715, 0.9999412634788912, 0.9999412634788912, 0.9999817919070126, 0.9999817919070126, 0.9999896830273372)
    27 = @(1.124442944769613, 0.9686981865604244, 0.9689954176457359, 1.124442944769613, 0.9689087623183813, 0.9686981865604244, 1.068374821632422, 0.9173777296412539, 0.9686981865604244, 0.9689954176457359, 1.046719181207674, 1.046719181207674, 1.053937728015924, 1.020712182991924, 1.020712182991924, 1.007708683884049, 1.007708683884049, 1.002799643761305)
    28 = @(1.033071397338449, 0.9918147868316274, 0.9916862453123849, 1.033071397338449, 0.9917237199623147, 0.9918147868316274, 1.01821901744236, 0.9778055018278902, 0.9918147868316274, 0.9916862453123849, 1.012378821325417, 1.012378821325417, 1.014325553364398, 1.00552414316082, 1.00552414316082, 1.002096804078522, 1.002096804078522, 1.000720111452504)
    29 = @(0.9972263944468244, 1.000675183502722, 1.000703540764506, 0.9972263944468244, 1.000695244007062, 1.000675183502722, 0.9984679695955438, 1.001881449358864, 1.000675183502722, 1.000703540764506, 0.9989649676056653, 0.9989649676056653, 0.9987993016022916, 0.9995350395713509, 0.999535039571351, 0.9998200755541938, 0.9998200755541938, 0.9999416302792539)
    30 = @(1.009223040097993, 0.9975645637998064, 0.9977662724989611, 1.009223040097993, 0.9977074643522817, 0.9975645637998064, 1.005026004071202, 0.9940818253572086, 0.9975645637998064, 0.9977662724989611, 1.003494656298477, 1.003494656298477, 1.004005105556052, 1.00151795879892, 1.00151795879892, 1.000529610049142, 1.000529610049142, 1.000228195029576)
    31 = @(1.01668392584168, 0.9956716285978403, 0.9959164536572774, 1.01668392584168, 0.9958450651231909, 0.9956716285978403, 1.009119525233704, 0.9891571757632716, 0.9956716285978403, 0.9959164536572774, 1.006300189749479, 1.006300189749479, 1.00723996824422, 1.002757336032266, 1.002757336032266, 1.000985909173659, 1.000985909173659, 1.000398962369494)
}
foreach ($r in $rowValues.Keys) {
    $vals = $rowValues[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, $i + 3).Value = $vals[$i]
    }
}

# ---- Step 5: dimension A1:T29 -> A1:T31 handled automatically by engine on save ----